# "use currency from account"
# Add a new "currencyCode" column to the accounts sheet (and a new "icon"
# column before the existing "sortOrder" column, matching the layout used
# on the other sheets), and make the accounts sheet the active tab.

$wb = $excel.ActiveWorkbook
$accounts = $wb.Worksheets.Item("accounts")

# The accounts header row currently ends at column L ("sortOrder"). Insert
# a new blank column there so "sortOrder" is pushed one column to the
# right, then fill the freed column with "icon" (matching the "banks" /
# "bills" sheets), and append a new "currencyCode" column at the end.
$accounts.Columns.Item(12).Insert()
$accounts.Cells.Item(1, 12).Value = "icon"
$accounts.Cells.Item(1, 14).Value = "currencyCode"

# Make "accounts" the active sheet/tab (it was "bills" before).
$accounts.Activate()
